$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Luka Doncic", "PG,SG", "Los Angeles Lakers"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Nickeil Alexander-Walker", "SG,SF", "Minnesota Timberwolves"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Isaiah Collier", "PG,SG", "Utah Jazz"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("De'Aaron Fox", "PG,SG", "San Antonio Spurs"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}
